# Re-create the "program" as described in the commit: within the
# "sections" worksheet, the GBP rows (previously grouped together at the
# bottom of the block) are redistributed so each insurance-period group
# (column E) has its own GBP line, and every group below shifts down by
# one row. Re-write rows 11-35 (columns E, O, Q, S) with the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sections")

$rows = @(
    @{ Row = 11; E = 2; O = "GBP"; Q = 43333333;  S = 23333333  },
    @{ Row = 12; E = 3; O = "USD"; Q = 50000000;  S = 65000000  },
    @{ Row = 13; E = 3; O = "CAD"; Q = 50000000;  S = 65000000  },
    @{ Row = 14; E = 3; O = "EUR"; Q = 50000000;  S = 65000000  },
    @{ Row = 15; E = 3; O = "AUD"; Q = 50000000;  S = 65000000  },
    @{ Row = 16; E = 3; O = "GBP"; Q = 33333333;  S = 43333333  },
    @{ Row = 17; E = 4; O = "USD"; Q = 100000000; S = 115000000 },
    @{ Row = 18; E = 4; O = "CAD"; Q = 100000000; S = 115000000 },
    @{ Row = 19; E = 4; O = "EUR"; Q = 100000000; S = 115000000 },
    @{ Row = 20; E = 4; O = "AUD"; Q = 100000000; S = 115000000 },
    @{ Row = 21; E = 4; O = "GBP"; Q = 66666666;  S = 76666666  },
    @{ Row = 22; E = 5; O = "USD"; Q = 100000000; S = 215000000 },
    @{ Row = 23; E = 5; O = "CAD"; Q = 100000000; S = 215000000 },
    @{ Row = 24; E = 5; O = "EUR"; Q = 100000000; S = 215000000 },
    @{ Row = 25; E = 5; O = "AUD"; Q = 100000000; S = 215000000 },
    @{ Row = 26; E = 5; O = "GBP"; Q = 66666666;  S = 143333333 },
    @{ Row = 27; E = 6; O = "USD"; Q = 100000000; S = 315000000 },
    @{ Row = 28; E = 6; O = "CAD"; Q = 100000000; S = 315000000 },
    @{ Row = 29; E = 6; O = "EUR"; Q = 100000000; S = 315000000 },
    @{ Row = 30; E = 6; O = "AUD"; Q = 100000000; S = 315000000 },
    @{ Row = 31; E = 6; O = "GBP"; Q = 66666666;  S = 210000000 },
    @{ Row = 32; E = 7; O = "USD"; Q = 150000000; S = 415000000 },
    @{ Row = 33; E = 7; O = "CAD"; Q = 150000000; S = 415000000 },
    @{ Row = 34; E = 7; O = "EUR"; Q = 150000000; S = 415000000 },
    @{ Row = 35; E = 7; O = "AUD"; Q = 150000000; S = 415000000 }
)

foreach ($r in $rows) {
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("Q" + $r.Row).Value = $r.Q
    $ws.Range("S" + $r.Row).Value = $r.S
}
